# Apply updated dSF (column F) values for the affected rows.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    4  = -7
    7  = -1
    9  = -3
    11 = -5
    12 = -6
    13 = -4
    15 = 4
    16 = -6
    22 = -3
    28 = 4
    29 = 0
    34 = 5
    38 = 4
    39 = 3
    40 = 1
    44 = 5
    46 = -1
    52 = -2
    60 = 4
    62 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
